$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value2 = 1163.1111
$ws.Range("J17").Value2 = 1537.0
$ws.Range("L17").Value2 = 4611.0
$ws.Range("N17").Value2 = -4947.0
$ws.Range("H112").Value2 = 3286.0
$ws.Range("J112").Value2 = 3206.9285
$ws.Range("L112").Value2 = 9620.7855
$ws.Range("N112").Value2 = -11836.7855
$ws.Range("H132").Value2 = 7261.727
$ws.Range("I132").Value2 = 8071.0
$ws.Range("J132").Value2 = 1394.5
$ws.Range("K132").Value2 = 24213.0
$ws.Range("L132").Value2 = 4183.5
$ws.Range("M132").Value2 = -21683.0
$ws.Range("N132").Value2 = -9243.5
$ws.Range("H136").Value2 = 0.0
$ws.Range("J136").Value2 = 0.0
$ws.Range("L136").Value2 = 0.0
$ws.Range("N136").ClearContents()
$ws.Range("H137").Value2 = 22733102.0
$ws.Range("I137").Value2 = 29413010.0
$ws.Range("K137").Value2 = 88239030.0
$ws.Range("M137").Value2 = -88236480.0
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value2 = 1922931.1
$ws.Range("I32").Value2 = 2418088.2
$ws.Range("K32").Value2 = 2418088.2
$ws.Range("M32").Value2 = -2417801.2
$ws.Range("H45").Value2 = 2724.75
$ws.Range("I45").Value2 = 1906.125
$ws.Range("K45").Value2 = 1906.125
$ws.Range("M45").Value2 = -1529.125
$ws.Range("H46").Value2 = 6671.143
$ws.Range("I46").Value2 = 3937.25
$ws.Range("K46").Value2 = 3937.25
$ws.Range("M46").Value2 = -3618.25
$ws.Range("H61").Value2 = 2944539.2
$ws.Range("I61").Value2 = 3223.4666
$ws.Range("K61").Value2 = 3223.4666
$ws.Range("M61").Value2 = -3011.4666
$ws.Range("H74").Value2 = 947172.2
$ws.Range("I74").Value2 = 1266176.1
$ws.Range("J74").Value2 = 11427.533
$ws.Range("K74").Value2 = 1266176.1
$ws.Range("L74").Value2 = 11427.533
$ws.Range("M74").Value2 = -1265302.1
$ws.Range("N74").Value2 = -13175.533
$ws.Range("H77").Value2 = 947172.2
$ws.Range("I77").Value2 = 1266176.1
$ws.Range("J77").Value2 = 11427.533
$ws.Range("K77").Value2 = 6330880.5
$ws.Range("L77").Value2 = 57137.66499999999
$ws.Range("M77").Value2 = -6326512.5
$ws.Range("N77").Value2 = -65873.665
$ws.Range("H122").Value2 = 1925.0
$ws.Range("I122").Value2 = 1803.8462
$ws.Range("K122").Value2 = 5411.5386
$ws.Range("M122").Value2 = -2961.5386
$ws.Range("H131").Value2 = 94999.5
$ws.Range("J131").Value2 = 94999.5
$ws.Range("L131").Value2 = 94999.5
$ws.Range("N131").Value2 = -105079.5
$ws.Range("H132").Value2 = 4820.492
$ws.Range("I132").Value2 = 3011.205
$ws.Range("J132").Value2 = 7534.423
$ws.Range("K132").Value2 = 9033.615
$ws.Range("L132").Value2 = 22603.269
$ws.Range("M132").Value2 = -6503.615
$ws.Range("N132").Value2 = -27663.269
$ws.Range("H136").Value2 = 2944539.2
$ws.Range("I136").Value2 = 3223.4666
$ws.Range("K136").Value2 = 9670.399800000001
$ws.Range("M136").Value2 = -7120.399800000001
$ws.Range("H139").Value2 = 106607.7
$ws.Range("J139").Value2 = 106607.7
$ws.Range("L139").Value2 = 106607.7
$ws.Range("N139").Value2 = -116887.7
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H13").Value2 = 70497.5
$ws.Range("J13").Value2 = 70497.5
$ws.Range("L13").Value2 = 70497.5
$ws.Range("N13").Value2 = -70833.5
$ws.Range("H134").Value2 = 5052970.0
$ws.Range("I134").Value2 = 2504.6667
$ws.Range("K134").Value2 = 7514.000100000001
$ws.Range("M134").Value2 = -4979.000100000001
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value2 = 5000.0
$ws.Range("I15").Value2 = 8000.0
$ws.Range("J15").Value2 = 2000.0
$ws.Range("K15").Value2 = 8000.0
$ws.Range("L15").Value2 = 2000.0
$ws.Range("M15").Value2 = -7830.0
$ws.Range("N15").Value2 = -2340.0
$ws.Range("H31").Value2 = 1549170.6
$ws.Range("I31").Value2 = 2194058.5
$ws.Range("K31").Value2 = 2194058.5
$ws.Range("M31").Value2 = -2193763.5
$ws.Range("H34").Value2 = 1549170.6
$ws.Range("I34").Value2 = 2194058.5
$ws.Range("K34").Value2 = 2194058.5
$ws.Range("M34").Value2 = -2193856.5
$ws.Range("H52").Value2 = 83890.0
$ws.Range("I52").Value2 = 0.0
$ws.Range("J52").Value2 = 83890.0
$ws.Range("K52").Value2 = 0.0
$ws.Range("L52").Value2 = 83890.0
$ws.Range("M52").ClearContents()
$ws.Range("N52").Value2 = -84478.0
$ws.Range("H99").Value2 = 10538.308
$ws.Range("I99").Value2 = 19607.416
$ws.Range("K99").Value2 = 19607.416
$ws.Range("M99").Value2 = -18109.416
$ws.Range("H105").Value2 = 18530.334
$ws.Range("I105").Value2 = 21236.4
$ws.Range("J105").Value2 = 5000.0
$ws.Range("K105").Value2 = 21236.4
$ws.Range("L105").Value2 = 5000.0
$ws.Range("M105").Value2 = -19489.4
$ws.Range("N105").Value2 = -8494.0
$ws.Range("H126").Value2 = 10538.308
$ws.Range("I126").Value2 = 19607.416
$ws.Range("K126").Value2 = 58822.24800000001
$ws.Range("M126").Value2 = -56352.24800000001
$ws.Range("H132").Value2 = 2981.5144
$ws.Range("I132").Value2 = 2833.9285
$ws.Range("J132").Value2 = 3571.8572
$ws.Range("K132").Value2 = 8501.7855
$ws.Range("L132").Value2 = 10715.5716
$ws.Range("M132").Value2 = -5971.7855
$ws.Range("N132").Value2 = -15775.5716
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value2 = 5952.3335
$ws.Range("I114").Value2 = 549.2857
$ws.Range("J114").Value2 = 7843.4
$ws.Range("K114").Value2 = 1647.8571
$ws.Range("L114").Value2 = 23530.2
$ws.Range("M114").Value2 = 1606.1429
$ws.Range("N114").Value2 = -30038.2
$ws.Range("H124").Value2 = 13459.75
$ws.Range("I124").Value2 = 10892.667
$ws.Range("K124").Value2 = 32678.001
$ws.Range("M124").Value2 = -27768.001
$ws.Range("H133").Value2 = 6753.2
$ws.Range("I133").Value2 = 4043.5386
$ws.Range("J133").Value2 = 11785.429
$ws.Range("K133").Value2 = 12130.6158
$ws.Range("L133").Value2 = 35356.287
$ws.Range("M133").Value2 = -7070.6158
$ws.Range("N133").Value2 = -45476.287
$ws.Range("H138").Value2 = 26509.5
$ws.Range("I138").Value2 = 58204.668
$ws.Range("J138").Value2 = 15944.444
$ws.Range("K138").Value2 = 174614.004
$ws.Range("L138").Value2 = 47833.33199999999
$ws.Range("M138").Value2 = -169474.004
$ws.Range("N138").Value2 = -58113.33199999999
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value2 = 11577.833
$ws.Range("I132").Value2 = 15958.286
$ws.Range("J132").Value2 = 5445.2
$ws.Range("K132").Value2 = 47874.858
$ws.Range("L132").Value2 = 16335.6
$ws.Range("M132").Value2 = -45344.858
$ws.Range("N132").Value2 = -21395.6
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H34").Value2 = 9943.625
$ws.Range("J34").Value2 = 12487.25
$ws.Range("L34").Value2 = 12487.25
$ws.Range("N34").Value2 = -12831.25
$ws.Range("H55").Value2 = 1247.2307
$ws.Range("I55").Value2 = 1056.6111
$ws.Range("K55").Value2 = 1056.6111
$ws.Range("M55").Value2 = -883.6111000000001
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value2 = 7248274.0
$ws.Range("I132").Value2 = 8773374.0
$ws.Range("K132").Value2 = 26320122.0
$ws.Range("M132").Value2 = -26317592.0
